$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "model_4_7_24",
    "model_4_7_22",
    "model_4_7_23",
    "model_4_7_21",
    "model_4_7_20",
    "model_4_7_19",
    "model_4_7_18",
    "model_4_7_17",
    "model_4_7_16",
    "model_4_7_15",
    "model_4_7_2",
    "model_4_7_14",
    "model_4_7_3",
    "model_4_7_13",
    "model_4_7_4",
    "model_4_7_12",
    "model_4_7_8",
    "model_4_7_5",
    "model_4_7_6",
    "model_4_7_11",
    "model_4_7_7",
    "model_4_7_10",
    "model_4_7_9",
    "model_4_7_1",
    "model_4_7_0"
)

$values = @(
    @(0.1468087291125882, -0.2579805597569595, -2.96553883763914, -1.333159986838776, 0.9442313313484192, 1.768062114715576, 4.138510704040527, 2.883567810058594),
    @(0.1504525720508422, -0.2159530789692503, -2.998993032830487, -1.321149150420607, 0.9401987195014954, 1.708993434906006, 4.173424243927002, 2.868723630905151),
    @(0.156851131379888, -0.1805375422529754, -3.006820293246807, -1.302939392434559, 0.9331173300743103, 1.659217715263367, 4.18159294128418, 2.84621787071228),
    @(0.1643413077167032, -0.1547171644497956, -2.98646469707262, -1.27930651400975, 0.9248279333114624, 1.622927784919739, 4.160349369049072, 2.817009925842285),
    @(0.1743969549673726, -0.134907839883194, -2.939406785451967, -1.248681145072607, 0.9136992692947388, 1.595086336135864, 4.111238956451416, 2.779159784317017),
    @(0.2059220490751215, -0.00129242807254859, -2.898293720765271, -1.15189921911758, 0.8788103461265564, 1.40729284286499, 4.068332195281982, 2.659546136856079),
    @(0.2811992826562635, 0.2333844646972482, -2.690891101704143, -0.9281965041953024, 0.795500636100769, 1.077460050582886, 3.851883172988892, 2.383070468902588),
    @(0.2968081189240965, 0.2500584402426042, -2.598305217236901, -0.8813653329034743, 0.7782262563705444, 1.054025173187256, 3.755258798599243, 2.325191736221313),
    @(0.3105010538131133, 0.2980048730872318, -2.567584182333376, -0.8402932749682983, 0.76307213306427, 0.9866376519203186, 3.723197937011719, 2.274430513381958),
    @(0.3424171681344276, 0.438094787108811, -2.541247298562319, -0.7454880544506748, 0.7277503609657288, 0.7897446155548096, 3.695712327957153, 2.157259941101074),
    @(0.3686903837790492, 0.3583658217399145, -2.018812740851146, -0.585886612868791, 0.698673665523529, 0.9018018245697021, 3.150489568710327, 1.960007667541504),
    @(0.3851637497659807, 0.6278470787549142, -2.509368682068669, -0.6185769565097274, 0.6804425120353699, 0.5230521559715271, 3.662442922592163, 2.000410079956055),
    @(0.4266905198288746, 0.6810375659238197, -2.213967813150333, -0.469172182457738, 0.6344846487045288, 0.4482942223548889, 3.354157209396362, 1.815759539604187),
    @(0.4430670784411042, 0.8689321972044195, -2.441490179791506, -0.4464613544958789, 0.6163606643676758, 0.18421271443367, 3.591603517532349, 1.787690997123718),
    @(0.4466501926137666, 0.7992680351951651, -2.290802108189048, -0.4285245560368078, 0.612395167350769, 0.2821240425109863, 3.434342861175537, 1.765522956848145),
    @(0.448176941773755, 0.8693608203795998, -2.403895195436971, -0.4312629872381906, 0.6107054948806763, 0.1836103051900864, 3.552368879318237, 1.768907308578491),
    @(0.4488697983893211, 0.8484742682093742, -2.363073755754094, -0.4276165697299232, 0.6099387407302856, 0.2129658609628677, 3.509767055511475, 1.764400720596313),
    @(0.4489745532063536, 0.8405570423657185, -2.34799205402754, -0.4263890652844082, 0.6098228096961975, 0.2240933328866959, 3.494027376174927, 1.76288366317749),
    @(0.4492421607617639, 0.8455690748880246, -2.355153041825931, -0.4262201278019213, 0.6095266342163086, 0.2170490324497223, 3.501500844955444, 1.762674808502197),
    @(0.4493516470640899, 0.8666343326713039, -2.390850030587242, -0.4277214742789379, 0.6094054579734802, 0.1874423027038574, 3.538754463195801, 1.764530301094055),
    @(0.4494317100521524, 0.8494886545428588, -2.359922268465385, -0.4257537250987604, 0.6093168258666992, 0.211540162563324, 3.506478071212769, 1.762098431587219),
    @(0.4500231481360416, 0.8633159750384704, -2.380318281982722, -0.4255349722069728, 0.6086622476577759, 0.1921061873435974, 3.527763605117798, 1.761828064918518),
    @(0.4503566302794608, 0.855085358459555, -2.363866806637135, -0.4239534298713545, 0.6082932353019714, 0.2036741375923157, 3.510594606399536, 1.759873390197754),
    @(0.4568268506300701, 0.651626020575937, -1.711910491130705, -0.287375474758401, 0.6011326909065247, 0.4896314442157745, 2.830200672149658, 1.591075778007507),
    @(0.84854599859772, 0.8893489180287196, 0.9230559864558447, 0.9028073392759569, 0.1676149517297745, 0.1555174887180328, 0.0803002119064331, 0.1201210469007492)
)

for ($i = 0; $i -lt 25; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    for ($j = 0; $j -lt 8; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$i][$j]
    }
}
